$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2327
$ws.Range("I40").Value = 1949.75
$ws.Range("J40").Value = 2830
$ws.Range("K40").Value = 1949.75
$ws.Range("L40").Value = 2830
$ws.Range("M40").Value = -1774.75
$ws.Range("N40").Value = -3180
$ws.Range("H116").Value = 3478.9473
$ws.Range("I116").Value = 3090.5
$ws.Range("J116").Value = 3910.5557
$ws.Range("K116").Value = 3090.5
$ws.Range("L116").Value = 3910.5557
$ws.Range("M116").Value = 351.5
$ws.Range("N116").Value = -10794.5557
$ws.Range("H137").Value = 1255.619
$ws.Range("I137").Value = 1253.0769
$ws.Range("J137").Value = 1259.75
$ws.Range("K137").Value = 3759.2307
$ws.Range("L137").Value = 3779.25
$ws.Range("M137").Value = -1209.2307
$ws.Range("N137").Value = -8879.25
$ws.Range("H138").Value = 1395.86
$ws.Range("I138").Value = 755.63635
$ws.Range("J138").Value = 1711.194
$ws.Range("K138").Value = 2266.90905
$ws.Range("L138").Value = 5133.582
$ws.Range("M138").Value = 2873.09095
$ws.Range("N138").Value = -15413.582
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H41").Value = 3482.5
$ws.Range("I41").Value = 2542.5715
$ws.Range("K41").Value = 2542.5715
$ws.Range("M41").Value = -2128.5715
$ws.Range("H110").Value = 2109.6
$ws.Range("I110").Value = 1608.5714
$ws.Range("J110").Value = 2548
$ws.Range("K110").Value = 1608.5714
$ws.Range("L110").Value = 2548
$ws.Range("M110").Value = 436.4286
$ws.Range("N110").Value = -6638
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 19231638
$ws.Range("I94").Value = 25000742
$ws.Range("J94").Value = 1293.3334
$ws.Range("K94").Value = 25000742
$ws.Range("L94").Value = 1293.3334
$ws.Range("M94").Value = -25000291
$ws.Range("N94").Value = -2195.3334
$ws.Range("M107").ClearContents()
$ws.Range("H107").Value = 2450
$ws.Range("I107").Value = 0
$ws.Range("K107").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 90910190
$ws.Range("I16").Value = 125000950
$ws.Range("J16").Value = 1492.6666
$ws.Range("K16").Value = 125000950
$ws.Range("L16").Value = 1492.6666
$ws.Range("M16").Value = -125000663
$ws.Range("N16").Value = -2066.6666
$ws.Range("H22").Value = 471.42856
$ws.Range("I22").Value = 399.75
$ws.Range("J22").Value = 567
$ws.Range("K22").Value = 399.75
$ws.Range("L22").Value = 567
$ws.Range("M22").Value = -49.75
$ws.Range("N22").Value = -1267
$ws.Range("H31").Value = 1708.0233
$ws.Range("I31").Value = 863.46155
$ws.Range("J31").Value = 2999.7058
$ws.Range("K31").Value = 863.46155
$ws.Range("L31").Value = 2999.7058
$ws.Range("M31").Value = -568.46155
$ws.Range("N31").Value = -3589.7058
$ws.Range("H34").Value = 1708.0233
$ws.Range("I34").Value = 863.46155
$ws.Range("J34").Value = 2999.7058
$ws.Range("K34").Value = 863.46155
$ws.Range("L34").Value = 2999.7058
$ws.Range("M34").Value = -661.46155
$ws.Range("N34").Value = -3403.7058
$ws.Range("H97").Value = 25000
$ws.Range("J97").Value = 25000
$ws.Range("L97").Value = 25000
$ws.Range("N97").Value = -26982
$ws.Range("H99").Value = 1596.76
$ws.Range("I99").Value = 1584.5
$ws.Range("K99").Value = 1584.5
$ws.Range("M99").Value = -86.5
$ws.Range("H107").Value = 647.3043
$ws.Range("I107").Value = 470.91666
$ws.Range("K107").Value = 470.91666
$ws.Range("M107").Value = 1449.08334
$ws.Range("H113").Value = 90910190
$ws.Range("I113").Value = 125000950
$ws.Range("J113").Value = 1492.6666
$ws.Range("K113").Value = 125000950
$ws.Range("L113").Value = 1492.6666
$ws.Range("M113").Value = -124998780
$ws.Range("N113").Value = -5832.6666
$ws.Range("H126").Value = 1596.76
$ws.Range("I126").Value = 1584.5
$ws.Range("K126").Value = 4753.5
$ws.Range("M126").Value = -2283.5
$ws.Range("H134").Value = 19608912
$ws.Range("I134").Value = 20834420
$ws.Range("J134").Value = 780
$ws.Range("K134").Value = 62503260
$ws.Range("L134").Value = 2340
$ws.Range("M134").Value = -62500725
$ws.Range("N134").Value = -7410
$ws.Range("H141").Value = 33860
$ws.Range("J141").Value = 33860
$ws.Range("L141").Value = 33860
$ws.Range("N141").Value = -44220
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 7935.643
$ws.Range("I56").Value = 7935.643
$ws.Range("K56").Value = 7935.643
$ws.Range("M56").Value = -7405.643
$ws.Range("H131").Value = 19232108
$ws.Range("I131").Value = 500000100
$ws.Range("J131").Value = 1388.7
$ws.Range("K131").Value = 1500000300
$ws.Range("L131").Value = 4166.1
$ws.Range("M131").Value = -1499995260
$ws.Range("N131").Value = -14246.1
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2985.375
$ws.Range("I113").Value = 1519
$ws.Range("K113").Value = 1519
$ws.Range("M113").Value = 651
$ws.Range("H132").Value = 2148.8
$ws.Range("I132").Value = 1710.3077
$ws.Range("K132").Value = 5130.9231
$ws.Range("M132").Value = -2600.9231
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2134.8572
$ws.Range("I7").Value = 2084.4
$ws.Range("J7").Value = 2261
$ws.Range("K7").Value = 2084.4
$ws.Range("L7").Value = 2261
$ws.Range("M7").Value = -1972.4
$ws.Range("N7").Value = -2485
$ws.Range("H40").Value = 6101
$ws.Range("I40").Value = 2250
$ws.Range("K40").Value = 2250
$ws.Range("M40").Value = -2114
$ws.Range("H46").Value = 11308.25
$ws.Range("I46").Value = 2243
$ws.Range("K46").Value = 2243
$ws.Range("M46").Value = -2055
$ws.Range("H126").Value = 2134.8572
$ws.Range("I126").Value = 2084.4
$ws.Range("J126").Value = 2261
$ws.Range("K126").Value = 6253.200000000001
$ws.Range("L126").Value = 6783
$ws.Range("M126").Value = -3783.200000000001
$ws.Range("N126").Value = -11723
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N31").ClearContents()
$ws.Range("H31").Value = 0
$ws.Range("J31").Value = 0
$ws.Range("L31").Value = 0
$ws.Range("H52").Value = 13294.5
$ws.Range("J52").Value = 13294.5
$ws.Range("L52").Value = 13294.5
$ws.Range("N52").Value = -13746.5
$ws.Range("H81").Value = 1099.3334
$ws.Range("I81").Value = 1149
$ws.Range("J81").Value = 1000
$ws.Range("K81").Value = 2298
$ws.Range("L81").Value = 2000
$ws.Range("M81").Value = -1237
$ws.Range("N81").Value = -4122
$ws.Range("H84").Value = 1099.3334
$ws.Range("I84").Value = 1149
$ws.Range("J84").Value = 1000
$ws.Range("K84").Value = 11490
$ws.Range("L84").Value = 10000
$ws.Range("M84").Value = -6186
$ws.Range("N84").Value = -20608
$ws.Range("H122").Value = 48463784
$ws.Range("I122").Value = 50402296
$ws.Range("K122").Value = 151206888
$ws.Range("M122").Value = -151204438
$ws.Range("H126").Value = 47619610
$ws.Range("I126").Value = 58824144
$ws.Range("J126").Value = 337.75
$ws.Range("K126").Value = 176472432
$ws.Range("L126").Value = 1013.25
$ws.Range("M126").Value = -176469962
$ws.Range("N126").Value = -5953.25
